$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: becomes MicroPizzaNetWithSE (SE-Ratio=8)
$ws.Range("A2").Value = "MicroPizzaNetWithSE (SE-Ratio=8)"
$ws.Range("B2").Value = 25
$ws.Range("C2").Value = 0.1
$ws.Range("D2").Value = 0.0625
$ws.Range("E2").Value = 0.25
$ws.Range("F2").Value = 662
$ws.Range("G2").Value = 2.84765625
$ws.Range("H2").Value = 0.7119140625
$ws.Range("I2").Value = 0.2221250534057617
$ws.Range("J2").Value = 7.77437686920166
$ws.Range("K2").Value = 102.3837890625
$ws.Range("L2").Value = 6.081669807434082

# Row 3: becomes MicroPizzaNet (Original)
$ws.Range("A3").Value = "MicroPizzaNet (Original)"
$ws.Range("B3").Value = 8.333333333333332
$ws.Range("C3").Value = 0.125
$ws.Range("D3").Value = 0.25
$ws.Range("E3").Value = 0.08333333333333333
$ws.Range("F3").Value = 582
$ws.Range("G3").Value = 2.53515625
$ws.Range("H3").Value = 0.6337890625
$ws.Range("I3").Value = 0.1268362998962402
$ws.Range("J3").Value = 4.439270496368408
$ws.Range("K3").Value = 102.0947265625
$ws.Range("L3").Value = 5.852807760238647

# Row 4: becomes MicroPizzaNetWithSE (SE-Ratio=4)
$ws.Range("A4").Value = "MicroPizzaNetWithSE (SE-Ratio=4)"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 742
$ws.Range("G4").Value = 3.16015625
$ws.Range("H4").Value = 0.7900390625
$ws.Range("I4").Value = 0.2347159385681152
$ws.Range("J4").Value = 8.215057849884035
$ws.Range("K4").Value = 102.4853515625
$ws.Range("L4").Value = 5.895447731018066
